$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CasesTab (B2) Cypher query: the "Cohort" lookup/column is being
# removed (now handled by the dedicated Cohort tab query instead), so strip
# the trailing OPTIONAL MATCH (co:cohort) return column from this query.
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in  ['T3N1M0', 'Not Applicable']  OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# Row 2 has one fewer wrapped line now that the Cohort column is gone, so
# shrink it to match rows 3/4 (which wrap to the same 290pt height).
$ws.Rows.Item(2).RowHeight = 290

# Reflect the author's saved cursor position/selection (A2 in view, B2 selected).
$ws.Range("B2").Select()
